# Insert three new bulleted list items after the paragraph that ends with
# "Transform Mapping Data Flow." in the "Aggregation: Transforms. ..." bullet.
#
# New items (same list: numId=3, ilvl=0, ind left=600 hanging=360):
#   - Aggregate Statements. Kinds Context.
#   - Aggregate Mappings. Statement Context.
#   - Aggregate Transforms. Mapping Context.

$d = $word.ActiveDocument

$anchorText = "Aggregation: Transforms. Aggregation Templates (Model / Domains). Assertions (matching input Statements). Generate / Match Mappings (apply) Generate / Match Statements (transforms). Aggregate Kinds Transforms. Aggregate SPO Transforms. Agregate Resources U Occurrences Transform. Apply: Occurrences Roles in Transforms Declaration (reified model entities): refer source Predicate role / type as Occurrence, Attribute, etc. in Matching results. Refer source Subject as destination Object (matching Kinds and Wrapper types). Match / Aggregate Context, Match / Aggregate Subject, Match / Aggregate Predicate, Match / Aggregate Object (functional contexts). Transform Mapping Data Flow."

$newItems = @(
    "Aggregate Statements. Kinds Context.",
    "Aggregate Mappings. Statement Context.",
    "Aggregate Transforms. Mapping Context."
)

$anchor = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`a", "`v")
    if ($text -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    foreach ($p in $d.Paragraphs) {
        $text = $p.Range.Text
        if ($text -like "*Transform Mapping Data Flow*") {
            $anchor = $p
            break
        }
    }
}

$prev = $anchor
foreach ($itemText in $newItems) {
    $prev.Range.InsertParagraphAfter()
    $newPara = $prev.Next()
    $newPara.Range.InsertAfter($itemText)
    $prev = $newPara
}
